$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '29.962.71'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c.Style = 'Normal'
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '1.921.38'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +1.37%  '
$c.Style = 'Normal'
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '
$c.Style = 'Normal'
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '320.29'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -1.49%  '
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'
$c = $ws.Range("D7")
$c.NumberFormat = '@'
$c.Value = '0.5067'
$c.Style = 'Normal'
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  -2.18%  '
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.62%  '
$c.Style = 'Normal'
$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.08341'
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  -0.73%  '
$c.Style = 'Normal'
$c = $ws.Range("B10")
$c.NumberFormat = '@'
$c.Value = 'OKB'
$c.Style = 'Normal'
$c = $ws.Range("C10")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c.Style = 'Normal'
$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '42.16'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -1.26%  '
$c.Style = 'Normal'
$c = $ws.Range("B11")
$c.NumberFormat = '@'
$c.Value = 'Polygon'
$c.Style = 'Normal'
$c = $ws.Range("C11")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c.Style = 'Normal'
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '1.105'
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -0.82%  '
$c.Style = 'Normal'
$c = $ws.Range("B12")
$c.NumberFormat = '@'
$c.Value = 'Solana'
$c.Style = 'Normal'
$c = $ws.Range("C12")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c.Style = 'Normal'
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '23.84'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  +2.78%  '
$c.Style = 'Normal'
$c = $ws.Range("B13")
$c.NumberFormat = '@'
$c.Value = 'WrappedEther'
$c.Style = 'Normal'
$c = $ws.Range("C13")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '1.914.98'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  +1.18%  '
$c.Style = 'Normal'
$c = $ws.Range("B14")
$c.NumberFormat = '@'
$c.Value = 'Polkadot'
$c.Style = 'Normal'
$c = $ws.Range("C14")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c.Style = 'Normal'
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '6.407'
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -0.51%  '
$c.Style = 'Normal'
$c = $ws.Range("B15")
$c.NumberFormat = '@'
$c.Value = 'Chainlink'
$c.Style = 'Normal'
$c = $ws.Range("C15")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c.Style = 'Normal'
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '7.243'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -1.10%  '
$c.Style = 'Normal'
$c = $ws.Range("B16")
$c.NumberFormat = '@'
$c.Value = 'BinanceUSD'
$c.Style = 'Normal'
$c = $ws.Range("C16")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c.Style = 'Normal'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -0.12%  '
$c.Style = 'Normal'
$c = $ws.Range("B17")
$c.NumberFormat = '@'
$c.Value = 'Litecoin'
$c.Style = 'Normal'
$c = $ws.Range("C17")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '92.39'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -1.89%  '
$c.Style = 'Normal'
$c = $ws.Range("B18")
$c.NumberFormat = '@'
$c.Value = 'ShibaInu'
$c.Style = 'Normal'
$c = $ws.Range("C18")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c.Style = 'Normal'
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '0.00001098'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  -1.09%  '
$c.Style = 'Normal'
$c = $ws.Range("B19")
$c.NumberFormat = '@'
$c.Value = 'TRON'
$c.Style = 'Normal'
$c = $ws.Range("C19")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c.Style = 'Normal'
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '0.06517'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -1.91%  '
$c.Style = 'Normal'
$c = $ws.Range("B20")
$c.NumberFormat = '@'
$c.Value = 'Avalanche'
$c.Style = 'Normal'
$c = $ws.Range("C20")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c.Style = 'Normal'
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '18.24'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  +0.14%  '
$c.Style = 'Normal'
$c = $ws.Range("B21")
$c.NumberFormat = '@'
$c.Value = 'Dai'
$c.Style = 'Normal'
$c = $ws.Range("C21")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c.Style = 'Normal'
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '1.002'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  +0.15%  '
$c.Style = 'Normal'
$c = $ws.Range("B22")
$c.NumberFormat = '@'
$c.Value = 'Uniswap'
$c.Style = 'Normal'
$c = $ws.Range("C22")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c.Style = 'Normal'
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '5.957'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c.Style = 'Normal'
$c = $ws.Range("B23")
$c.NumberFormat = '@'
$c.Value = 'WrappedBTC'
$c.Style = 'Normal'
$c = $ws.Range("C23")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c.Style = 'Normal'
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '29.972.83'
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c.Style = 'Normal'
$c = $ws.Range("B24")
$c.NumberFormat = '@'
$c.Value = 'Cosmos'
$c.Style = 'Normal'
$c = $ws.Range("C24")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c.Style = 'Normal'
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '11.33'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +0.36%  '
$c.Style = 'Normal'
$c = $ws.Range("B25")
$c.NumberFormat = '@'
$c.Value = 'Toncoin'
$c.Style = 'Normal'
$c = $ws.Range("C25")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c.Style = 'Normal'
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '2.194'
$c.Style = 'Normal'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -1.57%  '
$c.Style = 'Normal'
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '22.15'
$c.Style = 'Normal'
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +2.66%  '
$c.Style = 'Normal'
$c = $ws.Range("B27")
$c.NumberFormat = '@'
$c.Value = 'WrappedliquidstakedEther2.0'
$c.Style = 'Normal'
$c = $ws.Range("C27")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.Style = 'Normal'
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '2.133.29'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +1.05%  '
$c.Style = 'Normal'
$c = $ws.Range("B28")
$c.NumberFormat = '@'
$c.Value = 'Monero'
$c.Style = 'Normal'
$c = $ws.Range("C28")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.Style = 'Normal'
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '162.40'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  +0.39%  '
$c.Style = 'Normal'
$c = $ws.Range("B29")
$c.NumberFormat = '@'
$c.Value = 'LidoDAOToken'
$c.Style = 'Normal'
$c = $ws.Range("C29")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.Style = 'Normal'
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '2.324'
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  -1.23%  '
$c.Style = 'Normal'
$c = $ws.Range("B30")
$c.NumberFormat = '@'
$c.Value = 'BitcoinCash'
$c.Style = 'Normal'
$c = $ws.Range("C30")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c.Style = 'Normal'
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '129.20'
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.28%  '
$c.Style = 'Normal'
$c = $ws.Range("B31")
$c.NumberFormat = '@'
$c.Value = 'ImmutableX'
$c.Style = 'Normal'
$c = $ws.Range("C31")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.Style = 'Normal'
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '1.130'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  +3.54%  '
$c.Style = 'Normal'
$c = $ws.Range("B32")
$c.NumberFormat = '@'
$c.Value = 'Stellar'
$c.Style = 'Normal'
$c = $ws.Range("C32")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c.Style = 'Normal'
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '0.1038'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -1.67%  '
$c.Style = 'Normal'
$c = $ws.Range("B33")
$c.NumberFormat = '@'
$c.Value = 'Filecoin'
$c.Style = 'Normal'
$c = $ws.Range("C33")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c.Style = 'Normal'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '5.979'
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -1.88%  '
$c.Style = 'Normal'
$c = $ws.Range("B34")
$c.NumberFormat = '@'
$c.Value = 'HuobiToken'
$c.Style = 'Normal'
$c = $ws.Range("C34")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c.Style = 'Normal'
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '3.816'
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +1.94%  '
$c.Style = 'Normal'
$c = $ws.Range("B35")
$c.NumberFormat = '@'
$c.Value = 'VeChain'
$c.Style = 'Normal'
$c = $ws.Range("C35")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c.Style = 'Normal'
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '0.02456'
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -1.52%  '
$c.Style = 'Normal'
$c = $ws.Range("B36")
$c.NumberFormat = '@'
$c.Value = 'InternetComputer(DFINITY)'
$c.Style = 'Normal'
$c = $ws.Range("C36")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c.Style = 'Normal'
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '5.409'
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  +2.31%  '
$c.Style = 'Normal'
$c = $ws.Range("B37")
$c.NumberFormat = '@'
$c.Value = 'Hedera'
$c.Style = 'Normal'
$c = $ws.Range("C37")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c.Style = 'Normal'
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '0.06429'
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -1.82%  '
$c.Style = 'Normal'
$c = $ws.Range("B38")
$c.NumberFormat = '@'
$c.Value = 'Algorand'
$c.Style = 'Normal'
$c = $ws.Range("C38")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c.Style = 'Normal'
$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '0.2158'
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -1.72%  '
$c.Style = 'Normal'
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.6505'
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c.Style = 'Normal'
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '8.725'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -0.19%  '
$c.Style = 'Normal'
$c = $ws.Range("B41")
$c.NumberFormat = '@'
$c.Value = 'ARBITRUM'
$c.Style = 'Normal'
$c = $ws.Range("C41")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c.Style = 'Normal'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '1.196'
$c.Style = 'Normal'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  -2.20%  '
$c.Style = 'Normal'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '11.41'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -3.25%  '
$c.Style = 'Normal'
$c = $ws.Range("B43")
$c.NumberFormat = '@'
$c.Value = 'TrustWalletToken'
$c.Style = 'Normal'
$c = $ws.Range("C43")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c.Style = 'Normal'
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '1.220'
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -0.73%  '
$c.Style = 'Normal'
$c = $ws.Range("B44")
$c.NumberFormat = '@'
$c.Value = 'NEARProtocol'
$c.Style = 'Normal'
$c = $ws.Range("C44")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = 'Normal'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '2.231'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  +8.59%  '
$c.Style = 'Normal'
$c = $ws.Range("B45")
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'
$c.Style = 'Normal'
$c = $ws.Range("C45")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '13.35'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +1.11%  '
$c.Style = 'Normal'
$c = $ws.Range("B46")
$c.NumberFormat = '@'
$c.Value = 'Decentraland'
$c.Style = 'Normal'
$c = $ws.Range("C46")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c.Style = 'Normal'
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '0.6067'
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -0.51%  '
$c.Style = 'Normal'
$c = $ws.Range("B47")
$c.NumberFormat = '@'
$c.Value = 'PancakeSwap'
$c.Style = 'Normal'
$c = $ws.Range("C47")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c.Style = 'Normal'
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '3.641'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -1.52%  '
$c.Style = 'Normal'
$c = $ws.Range("B48")
$c.NumberFormat = '@'
$c.Value = 'EOS'
$c.Style = 'Normal'
$c = $ws.Range("C48")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$c.Style = 'Normal'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '1.210'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -2.05%  '
$c.Style = 'Normal'
$c = $ws.Range("B49")
$c.NumberFormat = '@'
$c.Value = 'Quant'
$c.Style = 'Normal'
$c = $ws.Range("C49")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c.Style = 'Normal'
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '122.28'
$c.Style = 'Normal'
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  -1.78%  '
$c.Style = 'Normal'
$c = $ws.Range("B50")
$c.NumberFormat = '@'
$c.Value = 'Aave'
$c.Style = 'Normal'
$c = $ws.Range("C50")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c.Style = 'Normal'
$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '78.97'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -0.10%  '
$c.Style = 'Normal'
$c = $ws.Range("B51")
$c.NumberFormat = '@'
$c.Value = 'WEMIXTOKEN'
$c.Style = 'Normal'
$c = $ws.Range("C51")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c.Style = 'Normal'
$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '1.129'
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -2.92%  '
$c.Style = 'Normal'
